$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 63 (2026-02-28): "科儀佈置+春酒" / "台南道場春酒中午11:00-13:00/不共修"
#   -> "科儀活動佈置" / "台南道場春酒中午11:00-13:00\n13:30開始進行元辰燈科儀佈置"
$ws.Range("C63").Value = "科儀活動佈置"
$ws.Range("D63").Value = "台南道場春酒中午11:00-13:00`n13:30開始進行元辰燈科儀佈置"
$ws.Range("D63").WrapText = $true
$ws.Rows(63).RowHeight = 34

# Row 64 (2026-03-01): "元辰燈科儀\n玉皇上帝天公祖祝壽" / (empty)
#   -> "科儀活動" / "台南道場115年元辰燈科儀暨玉皇上帝天公祖祝壽活動"
$ws.Range("C64").Value = "科儀活動"
$ws.Range("D64").Value = "台南道場115年元辰燈科儀暨玉皇上帝天公祖祝壽活動"
$ws.Rows(64).RowHeight = 16.4

# Row 65 (2026-08-30): "宮慶\n王母娘娘聖誔祝壽" / (empty)
#   -> "科儀活動" / "台南道場9週年宮慶暨王母娘娘聖誔祝壽活動"
$ws.Range("C65").Value = "科儀活動"
$ws.Range("D65").Value = "台南道場9週年宮慶暨王母娘娘聖誔祝壽活動"
$ws.Range("D65").WrapText = $true
$ws.Rows(65).RowHeight = 16.4

# Update the sheet view to match the saved cursor/scroll position.
$ws.Application.ActiveWindow.ScrollRow = 44
$ws.Range("D64").Select()
